$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "{firstname}"
$ws.Range("G1").Value = "{{expenses}}"
$ws.Range("B3").Value = "{{address}}"
$ws.Range("B5").Value = "{{hobbies}}"

$ws.Range("D7").Select()
